# previsao_retorno.xlsx - "atualizei dados bibi e add"
#
# The refreshed export bumps every "INATIVO - X meses sem comprar" status
# label (column J / situacao) by +0.1 months, and pulls fresh source figures
# for three client rows (55, 116, 117) whose underlying metrics/dates moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Bump every "INATIVO - X.X meses sem comprar" situacao label by 0.1.
#    Column J is "situacao" (column index 10). Row 1 is the header.
# ---------------------------------------------------------------------------
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 10)
    $val = $cell.Value2
    if ($val -ne $null -and $val -match '^INATIVO - ([\d\.]+) meses sem comprar$') {
        $num = [double]$matches[1]
        $newnum = $num + 0.1
        $newstr = "INATIVO - " + ("{0:N1}" -f $newnum) + " meses sem comprar"
        $cell.Value = $newstr
    }
}

# ---------------------------------------------------------------------------
# 2) Refreshed source rows with new probabilities / counts / dates.
# ---------------------------------------------------------------------------

# Row 55 - IGOR SOARES DOS ANJOS (id 5984)
$ws.Range("B55").Value = 0.83
$ws.Range("C55").Value = 0.83
$ws.Range("E55").Value = 41
$ws.Range("H55").Value = 45858.90815972222
$ws.Range("I55").Value = 45873.90815972222

# Row 116 - BEMOL S/A (id 28458)
$ws.Range("E116").Value = 16771
$ws.Range("H116").Value = 45856.73876157407
$ws.Range("I116").Value = 45857.73876157407

# Row 117 - RUY MENTA JUNIOR (id 28502)
$ws.Range("B117").Value = 0.25
$ws.Range("D117").Value = 0.33
$ws.Range("E117").Value = 11
$ws.Range("F117").Value = 0.33
$ws.Range("H117").Value = 45856.77859953704
$ws.Range("I117").Value = 45918.77859953704
